{"js": "// Replace the date line and the 25 \"three-digit x one-digit\" answer\n// expressions in the practice-sheet table with the next day's values.\n// Each old string is unique in the document, so a scoped, exact,\n// case-sensitive search-and-replace is safe for every entry.\nconst replacements = [\n  [\"2024-01-09 Tuesday\", \"2024-01-10 Wednesday\"],\n  [\"955\u00d75=4775\", \"425\u00d77=2975\"],\n  [\"451\u00d79=4059\", \"213\u00d77=1491\"],\n  [\"121\u00d76=726\", \"555\u00d75=2775\"],\n  [\"566\u00d73=1698\", \"956\u00d72=1912\"],\n  [\"375\u00d74=1500\", \"940\u00d75=4700\"],\n  [\"653\u00d73=1959\", \"766\u00d74=3064\"],\n  [\"650\u00d74=2600\", \"108\u00d77=756\"],\n  [\"927\u00d78=7416\", \"286\u00d78=2288\"],\n  [\"231\u00d73=693\", \"148\u00d79=1332\"],\n  [\"496\u00d79=4464\", \"651\u00d72=1302\"],\n  [\"544\u00d72=1088\", \"504\u00d77=3528\"],\n  [\"857\u00d78=6856\", \"472\u00d72=944\"],\n  [\"864\u00d72=1728\", \"988\u00d73=2964\"],\n  [\"390\u00d79=3510\", \"131\u00d75=655\"],\n  [\"743\u00d74=2972\", \"521\u00d76=3126\"],\n  [\"613\u00d75=3065\", \"798\u00d76=4788\"],\n  [\"854\u00d79=7686\", \"988\u00d76=5928\"],\n  [\"182\u00d76=1092\", \"625\u00d73=1875\"],\n  [\"383\u00d74=1532\", \"186\u00d76=1116\"],\n  [\"980\u00d76=5880\", \"404\u00d74=1616\"],\n  [\"978\u00d73=2934\", \"780\u00d77=5460\"],\n  [\"993\u00d76=5958\", \"903\u00d79=8127\"],\n  [\"733\u00d74=2932\", \"838\u00d79=7542\"],\n  [\"541\u00d76=3246\", \"259\u00d73=777\"],\n  [\"400\u00d78=3200\", \"287\u00d74=1148\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 \"three-digit x one-digit\" answer\n# expressions in the practice-sheet table to the next day's values.\n# Every old string is unique in the document, so a plain Find/Replace\n# (wdReplaceAll) scoped to the whole document body is safe for each pair.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-01-09 Tuesday\", \"2024-01-10 Wednesday\"),\n    @(\"955\u00d75=4775\", \"425\u00d77=2975\"),\n    @(\"451\u00d79=4059\", \"213\u00d77=1491\"),\n    @(\"121\u00d76=726\", \"555\u00d75=2775\"),\n    @(\"566\u00d73=1698\", \"956\u00d72=1912\"),\n    @(\"375\u00d74=1500\", \"940\u00d75=4700\"),\n    @(\"653\u00d73=1959\", \"766\u00d74=3064\"),\n    @(\"650\u00d74=2600\", \"108\u00d77=756\"),\n    @(\"927\u00d78=7416\", \"286\u00d78=2288\"),\n    @(\"231\u00d73=693\", \"148\u00d79=1332\"),\n    @(\"496\u00d79=4464\", \"651\u00d72=1302\"),\n    @(\"544\u00d72=1088\", \"504\u00d77=3528\"),\n    @(\"857\u00d78=6856\", \"472\u00d72=944\"),\n    @(\"864\u00d72=1728\", \"988\u00d73=2964\"),\n    @(\"390\u00d79=3510\", \"131\u00d75=655\"),\n    @(\"743\u00d74=2972\", \"521\u00d76=3126\"),\n    @(\"613\u00d75=3065\", \"798\u00d76=4788\"),\n    @(\"854\u00d79=7686\", \"988\u00d76=5928\"),\n    @(\"182\u00d76=1092\", \"625\u00d73=1875\"),\n    @(\"383\u00d74=1532\", \"186\u00d76=1116\"),\n    @(\"980\u00d76=5880\", \"404\u00d74=1616\"),\n    @(\"978\u00d73=2934\", \"780\u00d77=5460\"),\n    @(\"993\u00d76=5958\", \"903\u00d79=8127\"),\n    @(\"733\u00d74=2932\", \"838\u00d79=7542\"),\n    @(\"541\u00d76=3246\", \"259\u00d73=777\"),\n    @(\"400\u00d78=3200\", \"287\u00d74=1148\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
